$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to store a literal text value (matches the source
    # workbook's inlineStr cells) even when the text looks numeric, e.g.
    # "555.95" or "1.00", so Excel does not silently coerce it to a number.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "60.320.53"
$ws.Range("E2").Value = "  +4.15%  "
Set-TextValue "D3" "2.433.29"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "555.95"
$ws.Range("E5").Value = "  +2.27%  "
Set-TextValue "D6" "139.50"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("E9").Value = "  +5.00%  "
Set-TextValue "D10" "5.74"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  -2.03%  "
Set-TextValue "D13" "25.03"
$ws.Range("E13").Value = "  +5.53%  "
Set-TextValue "D14" "2.865.75"
$ws.Range("E14").Value = "  +3.23%  "
Set-TextValue "D15" "60.252.35"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("E16").Value = "  +4.07%  "
Set-TextValue "D17" "2.434.88"
$ws.Range("E17").Value = "  +3.41%  "
Set-TextValue "D18" "11.40"
$ws.Range("E18").Value = "  +6.05%  "
Set-TextValue "D20" "333.87"
$ws.Range("E20").Value = "  +0.91%  "
Set-TextValue "D21" "6.76"
$ws.Range("E21").Value = "  +0.84%  "
Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  -0.08%  "
Set-TextValue "D23" "65.32"
$ws.Range("E23").Value = "  +4.35%  "
$ws.Range("E24").Value = "  +3.68%  "
Set-TextValue "D25" "8.64"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  +7.16%  "
$ws.Range("E29").Value = "  +1.84%  "
$ws.Range("E30").Value = "  +3.32%  "
Set-TextValue "D31" "169.44"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +1.62%  "
Set-TextValue "D33" "18.78"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E35").Value = "  +6.13%  "
Set-TextValue "D36" "4.23"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  +11.36%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D40" "324.30"
$ws.Range("E40").Value = "  +12.33%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D41" "39.85"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("E42").Value = "  +1.46%  "
Set-TextValue "D43" "141.08"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D46" "19.64"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("B47").Value = "Polygon"
$ws.Range("C47").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D47" "0.416"
$ws.Range("E47").Value = "  +8.93%  "
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("E49").Value = "  +2.15%  "
Set-TextValue "D50" "17.94"
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("E51").Value = "  -0.17%  "
